$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the affected range to Text format so that numeric-looking strings
# (e.g. "599.64", "2.88") are stored as literal text, not converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '61.595.96'
$ws.Range("E2").Value = '  -1.41%  '
$ws.Range("D3").Value = '3.000.88'
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '599.64'
$ws.Range("E5").Value = '  +2.56%  '
$ws.Range("D6").Value = '144.77'
$ws.Range("E6").Value = '  -1.91%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = '2.998.91'
$ws.Range("E9").Value = '  -0.58%  '
$ws.Range("E10").Value = '  -1.14%  '
$ws.Range("D11").Value = '6.03'
$ws.Range("E11").Value = '  +5.52%  '
$ws.Range("D13").Value = '0.0000230'
$ws.Range("E13").Value = '  +0.12%  '
$ws.Range("D14").Value = '34.51'
$ws.Range("E14").Value = '  -0.70%  '
$ws.Range("E15").Value = '  +2.22%  '
$ws.Range("D16").Value = '3.496.21'
$ws.Range("E16").Value = '  -0.52%  '
$ws.Range("E17").Value = '  -1.05%  '
$ws.Range("D18").Value = '61.541.05'
$ws.Range("E18").Value = '  -1.43%  '
$ws.Range("D19").Value = '3.000.86'
$ws.Range("E19").Value = '  -0.44%  '
$ws.Range("D20").Value = '452.32'
$ws.Range("E20").Value = '  -1.90%  '
$ws.Range("D21").Value = '14.05'
$ws.Range("E21").Value = '  +1.01%  '
$ws.Range("D22").Value = '0.687'
$ws.Range("E22").Value = '  +0.28%  '
$ws.Range("D23").Value = '7.36'
$ws.Range("E23").Value = '  +0.43%  '
$ws.Range("D24").Value = '81.66'
$ws.Range("E24").Value = '  +2.09%  '
$ws.Range("E25").Value = '  -2.28%  '
$ws.Range("D26").Value = '10.81'
$ws.Range("E26").Value = '  +8.26%  '
$ws.Range("D27").Value = '11.99'
$ws.Range("E27").Value = '  -2.68%  '
$ws.Range("E28").Value = '  +0.16%  '
$ws.Range("E29").Value = '  +2.52%  '
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("D31").Value = '7.24'
$ws.Range("E31").Value = '  +1.19%  '
$ws.Range("E32").Value = '  -1.16%  '
$ws.Range("D33").Value = '27.51'
$ws.Range("E33").Value = '  +1.01%  '
$ws.Range("E34").Value = '  +3.19%  '
$ws.Range("D35").Value = '0.0₃0833'
$ws.Range("E35").Value = '  +5.76%  '
$ws.Range("E36").Value = '  -1.47%  '
$ws.Range("E37").Value = '  +0.64%  '
$ws.Range("D38").Value = '9.25'
$ws.Range("E38").Value = '  +2.32%  '
$ws.Range("D39").Value = '50.58'
$ws.Range("E39").Value = '  +0.34%  '
$ws.Range("D40").Value = '2.08'
$ws.Range("E40").Value = '  -1.73%  '
$ws.Range("E41").Value = '  +10.37%  '
$ws.Range("D42").Value = '2.88'
$ws.Range("D43").Value = '400.39'
$ws.Range("E43").Value = '  -4.09%  '
$ws.Range("D44").Value = '39.74'
$ws.Range("E44").Value = '  +4.68%  '
$ws.Range("E45").Value = '  +0.46%  '
$ws.Range("E46").Value = '  -1.35%  '
$ws.Range("D47").Value = '2.700.35'
$ws.Range("D48").Value = '131.77'
$ws.Range("E48").Value = '  +2.10%  '
$ws.Range("E49").Value = '  +0.15%  '
$ws.Range("E50").Value = '  -0.35%  '
$ws.Range("E51").Value = '  +1.94%  '

# Restore the original (default/Normal) cell style so no stray style index is
# left referenced on these cells, matching the source workbook formatting.
$ws.Range("D2:E51").Style = "Normal"

